$wb = $excel.ActiveWorkbook

# --- Sheets: rename + add "policies", fix tab order -----------------------
# Worksheets.Add() inserts the new sheet before the current first sheet, so
# after this call index 1 is the new (policies) sheet and index 2 is the
# original (warm_up) sheet.
$null = $wb.Worksheets.Add()
$wb.Worksheets.Item(1).Name = "policies"
$wb.Worksheets.Item(2).Name = "warm_up"

# Put warm_up back in front so the tab order is warm_up, policies.
$wb.Worksheets.Item(1).Move(2)

$wsWarm = $wb.Worksheets.Item(1)
$wsPol  = $wb.Worksheets.Item(2)

# --- warm_up: zero-out the second (Male/on_art) IPT initiation curve ------
$wsWarm.Range("D85").Value = 0

# --- warm_up: selection / active sheet state -------------------------------
$wsWarm.Range("D29").Select()

$wsPol.Activate()
$wsPol.Range("B2").Select()

# --- policies: populate example policy run outputs -------------------------
$wsPol.Range("A1").Value = "POLICY_ID"
$wsPol.Range("B1").Value = "G_SET"
$wsPol.Range("C1").Value = "on_art"
$wsPol.Range("D1").Value = "ipt_init_perc"

$policyRows = @(
    @(1, 1, "yes", 0.29),
    @(2, 1, "yes", 0.29),
    @(3, 1, "yes", 0.7),
    @(1, 1, "no",  0),
    @(2, 1, "no",  0),
    @(3, 1, "no",  0.7),
    @(1, 2, "yes", 0.27),
    @(2, 2, "yes", 0.27),
    @(3, 2, "yes", 0.75),
    @(1, 2, "no",  0.014),
    @(2, 2, "no",  0.014),
    @(3, 2, "no",  0.75)
)

$r = 2
foreach ($row in $policyRows) {
    $wsPol.Cells.Item($r, 1).Value = $row[0]
    $wsPol.Cells.Item($r, 2).Value = $row[1]
    $wsPol.Cells.Item($r, 3).Value = $row[2]
    $wsPol.Cells.Item($r, 4).Value = $row[3]
    $r++
}
